# edit.ps1 - apply the weekly CompStat crime-data refresh described in the commit
# "New crime data collected": bump the volume/number + report date range in the
# header, and refresh every Week-to-Date / 28-Day / Year-to-Date / 2-Year / 14-Year /
# 31-Year figure in the crime-complaints table (rows 14-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats - used below to carry over the correct number format/style
# whenever a cell switches between a text placeholder ("0" / "***.*") and a real number
# (or vice versa), so the written cell keeps matching the style of its row siblings.
$xlPasteFormats = -4122

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 31   Number  43"
$ws.Range("C9").Value = "Report Covering the Week  10/21/2024  Through  10/27/2024"

# --- Crime complaints table updates (rows 14-30) ---

# Row 14
$ws.Range("L14").Value = -66.666666666666
$ws.Range("M14").Value = -50
$ws.Range("N14").Value = -86.666666666666

# Row 15
$ws.Range("D15").Value = 1
$ws.Range("D16").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E15").Value = -100
$ws.Range("E16").Copy()
$ws.Range("E15").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("J15").Value = 18
$ws.Range("K15").Value = -16.666666666666
$ws.Range("M15").Value = 25
$ws.Range("N15").Value = -54.545454545454

# Row 16
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 300
$ws.Range("F16").Value = 12
$ws.Range("H16").Value = 20
$ws.Range("I16").Value = 107
$ws.Range("J16").Value = 157
$ws.Range("K16").Value = -31.847133757961
$ws.Range("L16").Value = -44.845360824742
$ws.Range("M16").Value = -33.540372670807
$ws.Range("N16").Value = -86.455696202531

# Row 17
$ws.Range("C17").Value = 4
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 90
$ws.Range("I17").Value = 168
$ws.Range("J17").Value = 188
$ws.Range("K17").Value = -10.638297872340
$ws.Range("L17").Value = -4
$ws.Range("M17").Value = 29.230769230769
$ws.Range("N17").Value = -64.705882352941

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -60
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 25
$ws.Range("H18").Value = -52
$ws.Range("I18").Value = 126
$ws.Range("J18").Value = 217
$ws.Range("K18").Value = -41.935483870967
$ws.Range("L18").Value = -61.111111111111
$ws.Range("M18").Value = -31.147540983606
$ws.Range("N18").Value = -83.486238532110

# Row 19
$ws.Range("C19").Value = 20
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = 9.090909090909
$ws.Range("I19").Value = 631
$ws.Range("J19").Value = 829
$ws.Range("K19").Value = -23.884197828709
$ws.Range("L19").Value = -26.456876456876
$ws.Range("M19").Value = 0
$ws.Range("N19").Value = -52.875280059746

# Row 20
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 250
$ws.Range("I20").Value = 37
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = -7.5
$ws.Range("M20").Value = -11.904761904761
$ws.Range("N20").Value = -91.991341991342

# Row 21
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 33.333333333333
$ws.Range("F21").Value = 123
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = 7.894736842105
$ws.Range("I21").Value = 1086
$ws.Range("J21").Value = 1447
$ws.Range("K21").Value = -24.948168624740
$ws.Range("L21").Value = -32.672039677619
$ws.Range("M21").Value = -6.620808254514
$ws.Range("N21").Value = -71.995874161939

# Row 22
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("D22").Copy()
$ws.Range("C22").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("M22").Value = -45

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = 37.5
$ws.Range("I23").Value = 75
$ws.Range("J23").Value = 108
$ws.Range("K23").Value = -30.555555555555
$ws.Range("L23").Value = -36.440677966101
$ws.Range("M23").Value = -25

# Row 24
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 22.727272727272
$ws.Range("F24").Value = 134
$ws.Range("G24").Value = 127
$ws.Range("H24").Value = 5.511811023622
$ws.Range("I24").Value = 1256
$ws.Range("J24").Value = 1237
$ws.Range("K24").Value = 1.535974130962
$ws.Range("L24").Value = -32.327586206896
$ws.Range("M24").Value = -11.111111111111

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 20
$ws.Range("F25").Value = 73
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = 14.0625
$ws.Range("I25").Value = 722
$ws.Range("J25").Value = 701
$ws.Range("K25").Value = 2.995720399429
$ws.Range("L25").Value = -43.149606299212

# Row 26
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = 25
$ws.Range("G26").Value = 31
$ws.Range("H26").Value = 45.161290322580
$ws.Range("I26").Value = 359
$ws.Range("J26").Value = 381
$ws.Range("K26").Value = -5.774278215223
$ws.Range("L26").Value = -7.948717948717
$ws.Range("M26").Value = -8.184143222506

# Row 27
$ws.Range("D27").Value = 2
$ws.Range("D26").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E27").Value = -100
$ws.Range("E26").Copy()
$ws.Range("E27").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("G27").Value = 2
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 31
$ws.Range("K27").Value = -41.935483870967
$ws.Range("L27").Value = -40

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("C26").Copy()
$ws.Range("C28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("D28").Value = 1
$ws.Range("D26").Copy()
$ws.Range("D28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("E28").Value = 100
$ws.Range("E26").Copy()
$ws.Range("E28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 58
$ws.Range("J28").Value = 47
$ws.Range("K28").Value = 23.404255319148
$ws.Range("L28").Value = -27.5

# Row 29
$ws.Range("L29").Value = -25

# Row 30
$ws.Range("L30").Value = -42.857142857142
